$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.866888523101807
$ws.Range("B1").Value = 2.351571083068848
$ws.Range("C1").Value = 4.395833015441895
$ws.Range("D1").Value = 2.920371294021606
$ws.Range("E1").Value = 0.6254712343215942
